$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 held the stale IP "192.168.1.113" - replace it with "127.0.0.1"
# (text format already applied to C2, so the new value stays a text shared string)
$ws.Range("C2").Value = "127.0.0.1"

# E2 ("192.168.0.24") now gets the same Text number format C2 uses
$ws.Range("E2").NumberFormat = "@"

# Move the active selection from C2 to E2
$ws.Range("E2").Select()
